$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 50: extend existing row (previously only had D50) with new data
# across columns D:I. D50 already holds the string that needs to stay
# (problem source), so just overwrite it with the correct text.
# ---------------------------------------------------------------------
$ws.Range("D50").Value = "运营反馈"
$ws.Range("E50").Value = "请求订单显示这个报错INTERNAL_SERVER_ERROR"
$ws.Range("F50").Value = "商户请求中的域名已经被删除"
$ws.Range("G50").Value = "运营的问题，向运营反馈问题运营，有运营和商户沟通解决方案"

# Dates for H50 / I50 - copy number format (style) from an existing
# date cell first so no new style/numFmt gets created, then set value.
$ws.Range("H46").Copy()
$ws.Range("H50").PasteSpecial(-4122)
$ws.Range("H50").Value = Get-Date -Year 2025 -Month 6 -Day 26 -Hour 15 -Minute 3 -Second 0

$ws.Range("H46").Copy()
$ws.Range("I50").PasteSpecial(-4122)
$ws.Range("I50").Value = Get-Date -Year 2025 -Month 6 -Day 26 -Hour 18 -Minute 7 -Second 0

# ---------------------------------------------------------------------
# Row 51: brand new row, columns D:F
# ---------------------------------------------------------------------
$ws.Range("D51").Value = "运营反馈"

# E51 must end up with no explicit cell style (relies on the column's
# default style), matching the plain-text cells elsewhere in column E
# (e.g. E46/E47). Copy formatting from such a cell first.
$ws.Range("E19").Copy()
$ws.Range("E51").PasteSpecial(-4122)
$ws.Range("E51").Value = "支付无法进入到订单页面，只有请求订单。并且该问题随机出现，随机解决"

# F51 needs the highlighted "status" style (same as the empty C column
# status cells), so copy that formatting over before setting the value.
$ws.Range("C46").Copy()
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("F51").Value = "排查中"

# ---------------------------------------------------------------------
# Update the view state to match the final workbook (scrolled one row
# further, selection moved to G55)
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G55").Select()

$excel.CutCopyMode = $false
